# Update "想去人数" (want-to-go count) figures across the workbook, and
# refresh the "本地生活" (Local Life) sheet: the "PLAVE with animate cafe"
# pop-up has ended and drops off the list, so every remaining local-life
# row shifts up by one and picks up its refreshed visitor count.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "展览" (Exhibitions) - F column "想去人数" updates
# ---------------------------------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value  = 29
$wsExpo.Range("F4").Value  = 1220
$wsExpo.Range("F5").Value  = 1661
$wsExpo.Range("F6").Value  = 888
$wsExpo.Range("F8").Value  = 2267
$wsExpo.Range("F9").Value  = 670
$wsExpo.Range("F10").Value = 550
$wsExpo.Range("F12").Value = 638
$wsExpo.Range("F14").Value = 163
$wsExpo.Range("F16").Value = 2087
$wsExpo.Range("F17").Value = 1218
$wsExpo.Range("F18").Value = 669
$wsExpo.Range("F21").Value = 19
$wsExpo.Range("F24").Value = 312
$wsExpo.Range("F27").Value = 521
$wsExpo.Range("F31").Value = 4482
$wsExpo.Range("F32").Value = 26

# ---------------------------------------------------------------------
# Sheet "演出" (Performances) - F column "想去人数" updates
# ---------------------------------------------------------------------
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F8").Value  = 56
$wsShow.Range("F24").Value = 20
$wsShow.Range("F25").Value = 234
$wsShow.Range("F34").Value = 59

# ---------------------------------------------------------------------
# Sheet "本地生活" (Local Life) - the "PLAVE with animate cafe" pop-up
# (row 6, id 5) has finished. Remove it: every later row moves up one
# slot, keeping its own id/"A" number, and the two rows that slide into
# place (the old NIJISANJI EN rows) carry their refreshed visitor counts.
# ---------------------------------------------------------------------
$wsLocal = $wb.Worksheets.Item("本地生活")

# Row 6 becomes what used to be row 7 (NIJISANJI EN official store),
# with its updated visitor count (482 -> 483) and now-numeric price.
# (Force text format first so the "yyyy-mm-dd"-looking date string isn't
# auto-converted into a real date value, then drop the format stamp again
# so the cell keeps using the sheet's default/general style.)
$wsLocal.Range("B6").NumberFormat = "@"
$wsLocal.Range("B6").Value = "2024-03-21"
$wsLocal.Range("B6").ClearFormats()
$wsLocal.Range("C6").Value = "上海·NIJISANJI EN 官方授权主题店"
$wsLocal.Range("D6").Value = "西藏北路166号（地铁8号线曲阜路下） 静安大悦城"
$wsLocal.Range("E6").Value = "2024.03.21 00:00-04.28 23:59"
$wsLocal.Range("F6").Value = 483
$wsLocal.Range("G6").Value = 30
$wsLocal.Range("H6").Value = "https://show.bilibili.com/platform/detail.html?id=82858"
$wsLocal.Range("I6").Value = "//i1.hdslb.com/bfs/openplatform/202403/eeG6Usri1710399065622.jpeg"

# Row 7 becomes what used to be row 8 (NIJISANJI EN x animate cafe),
# with its updated visitor count (126 -> 129).
$wsLocal.Range("B7").NumberFormat = "@"
$wsLocal.Range("B7").Value = "2024-04-24"
$wsLocal.Range("B7").ClearFormats()
$wsLocal.Range("C7").Value = "上海·「NIJISANJI EN x animate cafe」"
$wsLocal.Range("D7").Value = "西藏北路198号大悦城北座8楼N809-1 animate cafe上海店"
$wsLocal.Range("E7").Value = "2024.04.24 00:00-05.22 23:59"
$wsLocal.Range("F7").Value = 129
$wsLocal.Range("G7").Value = 30
$wsLocal.Range("H7").Value = "https://show.bilibili.com/platform/detail.html?id=83223"
$wsLocal.Range("I7").Value = "//i0.hdslb.com/bfs/openplatform/202403/LzJJK9lc1711096202393.jpeg"

# The old row 8 is now a duplicate of row 7's data; delete the trailing row.
$wsLocal.Rows.Item(8).Delete()

# F column updates for the rows above the ones that were rewritten.
$wsLocal.Range("F4").Value = 1398
$wsLocal.Range("F5").Value = 1777

# ---------------------------------------------------------------------
# Sheet "全部类型" (All Types) - F column "想去人数" updates
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value  = 1398
$wsAll.Range("F4").Value  = 1777
$wsAll.Range("F5").Value  = 483
$wsAll.Range("F9").Value  = 29
$wsAll.Range("F11").Value = 1220
$wsAll.Range("F12").Value = 1661
$wsAll.Range("F14").Value = 56
$wsAll.Range("F16").Value = 888
$wsAll.Range("F18").Value = 2267
$wsAll.Range("F19").Value = 670
$wsAll.Range("F20").Value = 550
$wsAll.Range("F22").Value = 638
$wsAll.Range("F25").Value = 163
$wsAll.Range("F29").Value = 2087
$wsAll.Range("F30").Value = 1218
$wsAll.Range("F31").Value = 669
$wsAll.Range("F37").Value = 19
$wsAll.Range("F39").Value = 129
$wsAll.Range("F42").Value = 20
$wsAll.Range("F43").Value = 521
$wsAll.Range("F47").Value = 4482
$wsAll.Range("F49").Value = 59
